$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1309.2858
$ws.Range("J17").Value = 1309.2858
$ws.Range("L17").Value = 3927.8574
$ws.Range("N17").Value = -4263.857400000001
$ws.Range("H18").Value = 1239.9
$ws.Range("I18").Value = 711
$ws.Range("K18").Value = 711
$ws.Range("M18").Value = -427
$ws.Range("H21").Value = 9179.25
$ws.Range("I21").Value = 717
$ws.Range("K21").Value = 717
$ws.Range("M21").Value = -249
$ws.Range("H23").Value = 9179.25
$ws.Range("I23").Value = 717
$ws.Range("K23").Value = 717
$ws.Range("M23").Value = -483
$ws.Range("H112").Value = 1008.0526
$ws.Range("J112").Value = 1074.2941
$ws.Range("L112").Value = 3222.8823
$ws.Range("N112").Value = -5438.8823
$ws.Range("H129").Value = 3306.1
$ws.Range("I129").Value = 6393.4707
$ws.Range("J129").Value = 1024.1305
$ws.Range("K129").Value = 19180.4121
$ws.Range("L129").Value = 3072.3915
$ws.Range("M129").Value = -14180.4121
$ws.Range("N129").Value = -13072.3915
$ws.Range("H136").Value = 58219
$ws.Range("J136").Value = 58219
$ws.Range("L136").Value = 58219
$ws.Range("N136").Value = -68419
$ws.Range("H138").Value = 2965.72
$ws.Range("I138").Value = 2378
$ws.Range("J138").Value = 3151.3157
$ws.Range("K138").Value = 7134
$ws.Range("L138").Value = 9453.947100000001
$ws.Range("M138").Value = -1994
$ws.Range("N138").Value = -19733.9471
$ws.Range("H139").Value = 69118.28999999999
$ws.Range("J139").Value = 69118.28999999999
$ws.Range("L139").Value = 69118.28999999999
$ws.Range("N139").Value = -79398.28999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11012.623
$ws.Range("I32").Value = 9726.145500000001
$ws.Range("J32").Value = 23362.8
$ws.Range("K32").Value = 9726.145500000001
$ws.Range("L32").Value = 23362.8
$ws.Range("M32").Value = -9439.145500000001
$ws.Range("N32").Value = -23936.8
$ws.Range("H61").Value = 2163.9167
$ws.Range("I61").Value = 1515.05
$ws.Range("K61").Value = 1515.05
$ws.Range("M61").Value = -1303.05
$ws.Range("H74").Value = 1385.2
$ws.Range("I74").Value = 1487.7778
$ws.Range("K74").Value = 1487.7778
$ws.Range("M74").Value = -613.7778000000001
$ws.Range("H77").Value = 1385.2
$ws.Range("I77").Value = 1487.7778
$ws.Range("K77").Value = 7438.889
$ws.Range("M77").Value = -3070.889
$ws.Range("H136").Value = 2163.9167
$ws.Range("I136").Value = 1515.05
$ws.Range("K136").Value = 4545.15
$ws.Range("M136").Value = -1995.15

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 36917.45
$ws.Range("I20").Value = 46081.434
$ws.Range("K20").Value = 46081.434
$ws.Range("M20").Value = -45834.434
$ws.Range("H99").Value = 2515.5454
$ws.Range("I99").Value = 3580
$ws.Range("J99").Value = 2409.1
$ws.Range("K99").Value = 3580
$ws.Range("L99").Value = 2409.1
$ws.Range("M99").Value = -2082
$ws.Range("N99").Value = -5405.1
$ws.Range("H138").Value = 60915
$ws.Range("J138").Value = 60915
$ws.Range("L138").Value = 60915
$ws.Range("N138").Value = -71195

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 117.111115
$ws.Range("I7").Value = 69.72727
$ws.Range("J7").Value = 191.57143
$ws.Range("K7").Value = 69.72727
$ws.Range("L7").Value = 191.57143
$ws.Range("M7").Value = 43.27273
$ws.Range("N7").Value = -417.57143
$ws.Range("H28").Value = 24423.834
$ws.Range("J28").Value = 24423.834
$ws.Range("L28").Value = 24423.834
$ws.Range("N28").Value = -24913.834
$ws.Range("H31").Value = 2901.6226
$ws.Range("I31").Value = 852.2857
$ws.Range("J31").Value = 3637.282
$ws.Range("K31").Value = 852.2857
$ws.Range("L31").Value = 3637.282
$ws.Range("M31").Value = -557.2857
$ws.Range("N31").Value = -4227.282
$ws.Range("H34").Value = 2901.6226
$ws.Range("I34").Value = 852.2857
$ws.Range("J34").Value = 3637.282
$ws.Range("K34").Value = 852.2857
$ws.Range("L34").Value = 3637.282
$ws.Range("M34").Value = -650.2857
$ws.Range("N34").Value = -4041.282
$ws.Range("H92").Value = 29999
$ws.Range("J92").Value = 29999
$ws.Range("L92").Value = 29999
$ws.Range("N92").Value = -34991
$ws.Range("H122").Value = 1201.6
$ws.Range("I122").Value = 1270
$ws.Range("J122").Value = 1099
$ws.Range("K122").Value = 3810
$ws.Range("L122").Value = 3297
$ws.Range("M122").Value = -1360
$ws.Range("N122").Value = -8197
$ws.Range("H134").Value = 1291.7059
$ws.Range("I134").Value = 804.53845
$ws.Range("K134").Value = 2413.61535
$ws.Range("M134").Value = 121.38465
$ws.Range("H140").Value = 54859.8
$ws.Range("J140").Value = 54859.8
$ws.Range("L140").Value = 54859.8
$ws.Range("N140").Value = -65219.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 483.92593
$ws.Range("J122").Value = 500.90475
$ws.Range("L122").Value = 4508.14275
$ws.Range("N122").Value = -9408.142749999999
$ws.Range("H131").Value = 830.01
$ws.Range("I131").Value = 600
$ws.Range("J131").Value = 834.7041
$ws.Range("K131").Value = 1800
$ws.Range("L131").Value = 2504.1123
$ws.Range("M131").Value = 3240
$ws.Range("N131").Value = -12584.1123
$ws.Range("H137").Value = 4046529.2
$ws.Range("I137").Value = 78953.08
$ws.Range("J137").Value = 8344737
$ws.Range("K137").Value = 236859.24
$ws.Range("L137").Value = 25034211
$ws.Range("M137").Value = -231759.24
$ws.Range("N137").Value = -25044411

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 123077.766
$ws.Range("I70").Value = 158092.47
$ws.Range("J70").Value = 9280
$ws.Range("K70").Value = 158092.47
$ws.Range("L70").Value = 9280
$ws.Range("M70").Value = -157822.47
$ws.Range("N70").Value = -9820
$ws.Range("H73").Value = 123077.766
$ws.Range("I73").Value = 158092.47
$ws.Range("J73").Value = 9280
$ws.Range("K73").Value = 158092.47
$ws.Range("L73").Value = 9280
$ws.Range("M73").Value = -157156.47
$ws.Range("N73").Value = -11152
$ws.Range("H113").Value = 1321.4
$ws.Range("I113").Value = 1069.5
$ws.Range("J113").Value = 1573.3
$ws.Range("K113").Value = 1069.5
$ws.Range("L113").Value = 1573.3
$ws.Range("M113").Value = 1100.5
$ws.Range("N113").Value = -5913.3
$ws.Range("H138").Value = 81714.28999999999
$ws.Range("J138").Value = 80666.664
$ws.Range("L138").Value = 80666.664
$ws.Range("N138").Value = -90946.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5687
$ws.Range("I132").Value = 6144.5557
$ws.Range("K132").Value = 18433.6671
$ws.Range("M132").Value = -15903.6671
